$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date header in column C (13-01-2023), formatted like B1 ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C1").Value = "13-01-2023"

# --- Full data block (rows 2-39), in the new order, with column C values added ---
# Columns: Label, ValorizadoAl (B, existing date), Valorizado new date (C)
$data = @(
    @('1810 Renta variable', 171714.82, 183982.94),
    @('1822 Raices Valores Negociables', 352711.1, 353916.93),
    @('Adcap IOL Acciones Argentina', 56920.39, 56640.22),
    @('Allaria Acciones', 67262.34, 67392.53),
    @('Alpha Acciones', 152915.14, 152779.73),
    @('Alpha Mega', 451910.52, 450369.4),
    @('Alpha planeam equil', 15716.18, 15765.09),
    @('Alpha renta balan global', 162096.04, 163183.04),
    @('Argenfunds', 13115.3, 13098.67),
    @('Balanz', 97946.41, 149791.54),
    @('Compass Crecimiento', 388942.39, 286403.95),
    @('Consultatio Renta Variable', 1393.15, 1251.93),
    @('Delta Acciones', 0, 0),
    @('Delta Select', 175020.68, 175316.94),
    @('FBA Acciones Argentinas', 199456.88, 199276.31),
    @('FBA Calificado', 200689.13, 200773.94),
    @('Fima Acciones', 320077.42, 320449.12),
    @('Fima PB Acciones', 419711.1, 438745.35),
    @('Goal Acciones Argentinas', 23193.57, 23136.87),
    @('Goal acciones plus', 5395.58, 5402.1),
    @('HF Acciones Argentinas', 119179.1, 119893.54),
    @('HF Acciones Lideres', 226105.8, 246902.99),
    @('IAM Renta Variable', 42422.76, 44342.64),
    @('IEB Value', 8769.280000000001, 8758.74),
    @('Lombardi', 41477.86, 41486.82),
    @('MAF', 8697.51, 8822.65),
    @('Megainver', 38466.38, 38601.49),
    @('Pellegrini Acciones', 134399.26, 134062.76),
    @('Pionero Acciones', 141478.88, 141372.46),
    @('Premier Renta Variable', 36804.22, 76993.42999999999),
    @('Quinquela Acciones', 123664.98, 123889.66),
    @('Rofex 20 Renta Variable', 87190.08, 86999.64999999999),
    @('Supefondo RV', 2522709.04, 2523074.61),
    @('Superfondo ', 4951585.53, 4953727.79),
    @('Supergestion', 788234.37, 787038.3100000001),
    @('Toronto Trust Multimercado', 54903.93, 54823.59),
    @('avg', 350063.25, 351346.33),
    @('total', 12602277.12, 12648467.73)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}
